# Update "ランサーズ" (Lancers) scrape-log sheet to the new snapshot taken at
# 2025-12-04 06:28:42. The table shrinks from 22 data rows (2-23) to 16 data rows
# (2-17): row 2 keeps its content (same top-ranked item), rows 3-17 absorb the
# next 15 items from the new scrape, and the previous rows 18-23 are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scrape data for rows 2..17, columns B..H (A/timestamp handled separately
# below since it is identical for every row).
$rows = @(
    @("大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5423720", 385, "🔥AI,Ai ◆効率化"),
    @("子ども向け音声AI先生|2秒以内応答PTT会話MVP", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446833", 310, "🔥AI,Ai"),
    @("製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419380", 298, "🔥AI,Ai"),
    @("【GAS×API】eBay販売管理スプレッドシートの自動化構築(初期1〜2ヶ月+月額保守)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5447137", 293, "🔥API ◆自動化 ◇管理"),
    @("【スマホアプリ開発】音声データ推定アプリの依頼", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446990", 175, "★スマホアプリ ◆開発 ◇アプリ"),
    @("自動出品ツールの開発", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446867", 128, "◆ツール,開発"),
    @("製造業向けMR業務支援アプリケーションの機能開発エンジニア募集(Unity/C#)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5441557", 93, "◆開発 ◇アプリ"),
    @("初回 【Unity/常設】画像差替可能な「お絵かきシステム」開発依頼(多展開前提・2月末納品)", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5447021", 90, "◆開発"),
    @("製造業向け 技能習得・作業トレーニングVRシステムの開発(Unity/R3)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5441568", 83, "◆開発"),
    @("【時期未定】PERLからPHPへの生産管理システム移行依頼", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446668", 80, "◇管理 ○PHP"),
    @("【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5431107", 68, "◆開発"),
    @("【急募】MT5当人カスタムのパネル開発の依頼探してます。", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446849", 68, "◆開発"),
    @("【急募】宝くじ予想サイトの構築", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446997", 38, "◇サイト"),
    @("【急募】住所リストから太陽光パネル設置を自動判定するシステム", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5447102", 33, $null),
    @("【SESエンジニア募集】多様なプロジェクトに参画可能!", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5437544", 25, $null),
    @("注目 【電子工作】蓋を開くとmp3再生するスピーカー制作依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5446806", 13, $null),
)

# 1) Drop the 6 rows that no longer exist in the new snapshot (18-23), which
#    also shrinks the sheet dimension from H23 to H17.
$ws.Rows("18:23").Delete()

# 2) Column width tweaks (B 44->49, D 30->32, H 17->18). Excel stores column
#    widths in "characters", but round-trips whatever is assigned through a
#    pixel-based quantization step (~ +0.833 for this font); back that
#    constant out so the persisted width lands exactly on the target integer.
$ws.Columns.Item(2).ColumnWidth = 49 - 1 + 1/12
$ws.Columns.Item(4).ColumnWidth = 32 - 1 + 1/12
$ws.Columns.Item(8).ColumnWidth = 18 - 1 + 1/12

# 3) Refresh the "取得日時" (fetched-at) timestamp for every remaining data row.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-12-04 06:28:42"
}

# 4) Write the new title/category/price/deadline/url/score/skill-summary
#    values row by row.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    if ($row[6] -eq $null) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $row[6]
    }
}

# 5) Hyperlinks: the engine's Hyperlinks.Delete() clears the whole sheet
#    collection regardless of which range/item it was fetched from, so rebuild
#    the F-column links from scratch for the surviving rows rather than trying
#    to remove just the stale F18:F23 entries in place.
$ws.Hyperlinks.Delete()
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[4])
}
